# Insert a new "建物" (Building) worksheet between "土地" (Land) and
# "債務" (Debt), populated with a single building record, mirroring the
# structure/styling of the "土地" sheet.

$wb = $excel.ActiveWorkbook

$landSheet = $wb.Worksheets.Item(1)

# Duplicate the "土地" sheet (carries over column layout + cell styles)
# and drop it immediately after it; rename to "建物".
$landSheet.Copy($null, $landSheet)
$buildingSheet = $wb.Worksheets.Item(2)
$buildingSheet.Name = "建物"

# The source sheet had two data rows; the building sheet only needs one.
$buildingSheet.Rows.Item(3).Delete()

# Overwrite the remaining data row with the building record.
$buildingSheet.Range("A2").Value = 19
$buildingSheet.Range("B2").Value = "臺中市沙鹿區屏西路"
$buildingSheet.Range("C2").Value = 432
$buildingSheet.Range("D2").Value = "全部"
$buildingSheet.Range("E2").Value = "顏清標"
$buildingSheet.Range("F2").Value = "83年10月14日"
$buildingSheet.Range("G2").Value = "買賣"
$buildingSheet.Range("H2").Value = "(超過五年）"
$buildingSheet.Range("I2").Value = "building"
$buildingSheet.Range("J2").Value = "normal"

# "2012-11-28" looks like an ISO date -- force the cell to text first so
# Excel stores it as a literal string instead of reinterpreting it as a
# date serial, then reset the style to a plain (non-custom-format) style
# that matches its row neighbours.
$buildingSheet.Range("K2").NumberFormat = "@"
$buildingSheet.Range("K2").Value = "2012-11-28"
$buildingSheet.Range("K2").Style = $buildingSheet.Range("L2").Style

$buildingSheet.Range("L2").Value = "顏清標"
$buildingSheet.Range("M2").Value = 979
$buildingSheet.Range("N2").Value = "tmp68961"
$buildingSheet.Range("O2").Value = 19
$buildingSheet.Range("P2").Value = 1
$buildingSheet.Range("Q2").Value = 432

# Restore "土地" as the selected/active sheet (copying moved focus to
# the new sheet).
$landSheet.Select()

Write-Output "inserted 建物 sheet"
